$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.069.40'
$ws.Range('E2').Value = '  -1.38%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.789.70'
$ws.Range('E3').Value = '  -2.14%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.010'
$ws.Range('E4').Value = '  -0.88%  '
$ws.Range('B5').Value = 'USDC'
$ws.Range('C5').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.005'
$ws.Range('E5').Value = '  -1.11%  '
$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '307.39'
$ws.Range('E6').Value = '  -2.27%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4162'
$ws.Range('E7').Value = '  -1.61%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3524'
$ws.Range('E8').Value = '  -3.80%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07014'
$ws.Range('E9').Value = '  -3.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8382'
$ws.Range('E10').Value = '  -2.13%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.96'
$ws.Range('E11').Value = '  -2.49%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.866.74'
$ws.Range('E12').Value = '  -0.03%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.256'
$ws.Range('E13').Value = '  -1.72%  '
$ws.Range('B14').Value = 'TRON'
$ws.Range('C14').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.06833'
$ws.Range('E14').Value = '  -2.15%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.288'
$ws.Range('E15').Value = '  -2.89%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.009'
$ws.Range('E16').Value = '  -1.13%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '79.42'
$ws.Range('E17').Value = '  -0.21%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008665'
$ws.Range('E18').Value = '  -2.88%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.005'
$ws.Range('E19').Value = '  +0.32%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.96'
$ws.Range('E20').Value = '  -1.26%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '27.344.58'
$ws.Range('E21').Value = '  -1.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.022'
$ws.Range('E22').Value = '  +0.65%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.63'
$ws.Range('E23').Value = '  -0.23%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.038.32'
$ws.Range('E24').Value = '  -3.65%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.959'
$ws.Range('E25').Value = '  -0.96%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '153.14'
$ws.Range('E26').Value = '  -0.29%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.98'
$ws.Range('E27').Value = '  -1.94%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.000'
$ws.Range('E28').Value = '  -3.66%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '111.53'
$ws.Range('E29').Value = '  -3.06%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.648'
$ws.Range('E30').Value = '  -9.45%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08874'
$ws.Range('E31').Value = '  +0.09%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.937'
$ws.Range('E32').Value = '  -1.62%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7180'
$ws.Range('E33').Value = '  -5.47%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.321'
$ws.Range('E34').Value = '  -4.20%  '
$ws.Range('B35').Value = 'Frax'
$ws.Range('C35').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.004'
$ws.Range('E35').Value = '  -1.10%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.064'
$ws.Range('E36').Value = '  -4.93%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.067'
$ws.Range('E37').Value = '  -3.92%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01889'
$ws.Range('E38').Value = '  -2.65%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05076'
$ws.Range('E39').Value = '  -4.64%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.690'
$ws.Range('E40').Value = '  -5.12%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1609'
$ws.Range('E41').Value = '  -2.41%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4889'
$ws.Range('E42').Value = '  -2.40%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.223'
$ws.Range('E43').Value = '  -7.72%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.977'
$ws.Range('E44').Value = '  -4.13%  '
$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.004'
$ws.Range('E45').Value = '  -1.22%  '
$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '104.03'
$ws.Range('E46').Value = '  +0.19%  '
$ws.Range('E47').Value = '  -2.12%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.06325'
$ws.Range('E48').Value = '  -3.22%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4507'
$ws.Range('E49').Value = '  -2.57%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.574'
$ws.Range('E50').Value = '  -1.68%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '61.84'
$ws.Range('E51').Value = '  -2.41%  '
